$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 34673.668
$ws.Range("I20").Value = 34673.668
$ws.Range("K20").Value = 34673.668
$ws.Range("M20").Value = -34443.668
$ws.Range("H35").Value = 34673.668
$ws.Range("I35").Value = 34673.668
$ws.Range("K35").Value = 34673.668
$ws.Range("M35").Value = -34294.668
$ws.Range("H40").Value = 2060
$ws.Range("I40").Value = 1575
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 1575
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -1400
$ws.Range("N40").Value = -4350
$ws.Range("H112").Value = 10064.714
$ws.Range("J112").Value = 13740.6
$ws.Range("L112").Value = 41221.8
$ws.Range("N112").Value = -43437.8
$ws.Range("H129").Value = 1339.125
$ws.Range("I129").Value = 740
$ws.Range("J129").Value = 1698.6
$ws.Range("K129").Value = 2220
$ws.Range("L129").Value = 5095.799999999999
$ws.Range("M129").Value = 2780
$ws.Range("N129").Value = -15095.8
$ws.Range("H132").Value = 1218.82
$ws.Range("I132").Value = 1016.13336
$ws.Range("J132").Value = 3043
$ws.Range("K132").Value = 3048.40008
$ws.Range("L132").Value = 9129
$ws.Range("M132").Value = -518.4000800000003
$ws.Range("N132").Value = -14189
$ws.Range("H134").Value = 111770.24
$ws.Range("J134").Value = 111770.24
$ws.Range("L134").Value = 111770.24
$ws.Range("N134").Value = -121910.24
$ws.Range("H138").Value = 3962.94
$ws.Range("I138").Value = 3145.0425
$ws.Range("J138").Value = 4688.245
$ws.Range("K138").Value = 9435.127500000001
$ws.Range("L138").Value = 14064.735
$ws.Range("M138").Value = -4295.127500000001
$ws.Range("N138").Value = -24344.735

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32628.174
$ws.Range("I32").Value = 25782.904
$ws.Range("K32").Value = 25782.904
$ws.Range("M32").Value = -25495.904
$ws.Range("H61").Value = 2235.8293
$ws.Range("I61").Value = 2000.931
$ws.Range("J61").Value = 2803.5
$ws.Range("K61").Value = 2000.931
$ws.Range("L61").Value = 2803.5
$ws.Range("M61").Value = -1788.931
$ws.Range("N61").Value = -3227.5
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
$ws.Range("H136").Value = 2235.8293
$ws.Range("I136").Value = 2000.931
$ws.Range("J136").Value = 2803.5
$ws.Range("K136").Value = 6002.793
$ws.Range("L136").Value = 8410.5
$ws.Range("M136").Value = -3452.793
$ws.Range("N136").Value = -13510.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 29663.334
$ws.Range("J62").Value = 29663.334
$ws.Range("L62").Value = 29663.334
$ws.Range("N62").Value = -31035.334
$ws.Range("H65").Value = 29663.334
$ws.Range("J65").Value = 29663.334
$ws.Range("L65").Value = 88990.00199999999
$ws.Range("N65").Value = -95854.00199999999
$ws.Range("H75").Value = 64027.2
$ws.Range("I75").Value = 9309.817999999999
$ws.Range("J75").Value = 130904
$ws.Range("K75").Value = 9309.817999999999
$ws.Range("L75").Value = 130904
$ws.Range("M75").Value = -8373.817999999999
$ws.Range("N75").Value = -132776
$ws.Range("H78").Value = 64027.2
$ws.Range("I78").Value = 9309.817999999999
$ws.Range("J78").Value = 130904
$ws.Range("K78").Value = 27929.454
$ws.Range("L78").Value = 392712
$ws.Range("M78").Value = -23249.454
$ws.Range("N78").Value = -402072
$ws.Range("H105").Value = 2884.389
$ws.Range("I105").Value = 2744.9375
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 2744.9375
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -997.9375
$ws.Range("N105").Value = -7494
$ws.Range("H134").Value = 2472.5322
$ws.Range("I134").Value = 1742.7858
$ws.Range("J134").Value = 4005
$ws.Range("K134").Value = 5228.357400000001
$ws.Range("L134").Value = 12015
$ws.Range("M134").Value = -2693.357400000001
$ws.Range("N134").Value = -17085

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5616.5
$ws.Range("I31").Value = 4519.4287
$ws.Range("J31").Value = 6469.778
$ws.Range("K31").Value = 4519.4287
$ws.Range("L31").Value = 6469.778
$ws.Range("M31").Value = -4224.4287
$ws.Range("N31").Value = -7059.778
$ws.Range("H34").Value = 5616.5
$ws.Range("I34").Value = 4519.4287
$ws.Range("J34").Value = 6469.778
$ws.Range("K34").Value = 4519.4287
$ws.Range("L34").Value = 6469.778
$ws.Range("M34").Value = -4317.4287
$ws.Range("N34").Value = -6873.778
$ws.Range("H58").Value = 1289.9
$ws.Range("I58").Value = 1395.9032
$ws.Range("K58").Value = 1395.9032
$ws.Range("M58").Value = -1192.9032
$ws.Range("H97").Value = 25098.5
$ws.Range("I97").Value = 10000
$ws.Range("J97").Value = 40197
$ws.Range("K97").Value = 10000
$ws.Range("L97").Value = 40197
$ws.Range("N97").Value = -42179
$ws.Range("M97").Value = -9009
$ws.Range("H122").Value = 2156.25
$ws.Range("I122").Value = 2364.4
$ws.Range("J122").Value = 1809.3334
$ws.Range("K122").Value = 7093.200000000001
$ws.Range("L122").Value = 5428.0002
$ws.Range("M122").Value = -4643.200000000001
$ws.Range("N122").Value = -10328.0002
$ws.Range("H132").Value = 1583.5172
$ws.Range("I132").Value = 1500.3462
$ws.Range("J132").Value = 2304.3333
$ws.Range("K132").Value = 4501.0386
$ws.Range("L132").Value = 6912.999899999999
$ws.Range("M132").Value = -1971.0386
$ws.Range("N132").Value = -11972.9999
$ws.Range("H134").Value = 2118.7727
$ws.Range("I134").Value = 2117.6667
$ws.Range("J134").Value = 2120.1
$ws.Range("K134").Value = 6353.000100000001
$ws.Range("L134").Value = 6360.299999999999
$ws.Range("M134").Value = -3818.000100000001
$ws.Range("N134").Value = -11430.3
$ws.Range("H136").Value = 1289.9
$ws.Range("I136").Value = 1395.9032
$ws.Range("K136").Value = 4187.7096
$ws.Range("M136").Value = -1637.7096

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 2650
$ws.Range("J54").Value = 2650
$ws.Range("L54").Value = 7950
$ws.Range("N54").Value = -9068
$ws.Range("H55").Value = 3199.8823
$ws.Range("J55").Value = 3199.8823
$ws.Range("L55").Value = 9599.6469
$ws.Range("N55").Value = -9953.6469
$ws.Range("H86").Value = 734.3333
$ws.Range("J86").Value = 734.3333
$ws.Range("L86").Value = 2202.9999
$ws.Range("N86").Value = -4574.9999
$ws.Range("H89").Value = 734.3333
$ws.Range("J89").Value = 734.3333
$ws.Range("L89").Value = 6608.9997
$ws.Range("N89").Value = -18464.9997
$ws.Range("H131").Value = 15154356
$ws.Range("I131").Value = 17235
$ws.Range("J131").Value = 16668069
$ws.Range("K131").Value = 51705
$ws.Range("L131").Value = 50004207
$ws.Range("M131").Value = -46665
$ws.Range("N131").Value = -50014287

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 50000
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H102").Value = 1692.5151
$ws.Range("I102").Value = 1619.6666
$ws.Range("J102").Value = 2020.3334
$ws.Range("K102").Value = 1619.6666
$ws.Range("L102").Value = 2020.3334
$ws.Range("M102").Value = 2.333399999999983
$ws.Range("N102").Value = -5264.3334
$ws.Range("H123").Value = 33303.535
$ws.Range("J123").Value = 33303.535
$ws.Range("L123").Value = 33303.535
$ws.Range("N123").Value = -38203.535

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 12154806
$ws.Range("I122").Value = 14712028
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 44136084
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -44133634
$ws.Range("N122").Value = -28900

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 27779544
$ws.Range("I122").Value = 41667770
$ws.Range("J122").Value = 3101.6667
$ws.Range("K122").Value = 125003310
$ws.Range("L122").Value = 9305.000100000001
$ws.Range("M122").Value = -125000860
$ws.Range("N122").Value = -14205.0001
$ws.Range("H123").Value = 23613.5
$ws.Range("J123").Value = 23613.5
$ws.Range("L123").Value = 23613.5
$ws.Range("N123").Value = -33413.5
